$p = $ppt.ActivePresentation

# The deck ships two themes: theme1.xml ("Default" palette, wired to the
# notes master) and theme2.xml ("Simple Light" palette, wired to the slide
# master / active design). The commit swaps which palette lives in which
# theme part - the slide master's design becomes the "Default" colors
# while the notes master's design becomes the former "Simple Light"
# colors.
#
# Drive this through the exposed Theme/ThemeColorScheme COM surface: walk
# every theme color slot (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# on the active design and push in the "Default" palette's RGB values,
# matching the clrScheme that theme1.xml carried before the edit.

$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# index : RGB (decimal, 0xBBGGRR order expected by the RGB property)
$colorScheme.Item(1).RGB  = 0         # dk1      000000
$colorScheme.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 5800213   # dk2      158158
$colorScheme.Item(4).RGB  = 15987699  # lt2      F3F3F3
$colorScheme.Item(5).RGB  = 13077765  # accent1  058DC7
$colorScheme.Item(6).RGB  = 3322960   # accent2  50B432
$colorScheme.Item(7).RGB  = 1791725   # accent3  ED561B
$colorScheme.Item(8).RGB  = 61421     # accent4  EDEF00
$colorScheme.Item(9).RGB  = 15059748  # accent5  24CBE5
$colorScheme.Item(10).RGB = 7529828   # accent6  64E572
$colorScheme.Item(11).RGB = 13369378  # hlink    2200CC
$colorScheme.Item(12).RGB = 9116245   # folHlink 551A8B
